$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Memory Usage (bytes)"

$ws.Range("C2").Value = 18.01490783691406
$ws.Range("C3").Value = 16.25180244445801
$ws.Range("C4").Value = 16.78109169006348
$ws.Range("C5").Value = 16.17884635925293
$ws.Range("C6").Value = 15.8839225769043
